# debug(wrapper): Disable infor dict modification to test for memory leak
#
# On the DQN_Runs log sheet, the run-config row for the active run (row 2)
# had its time-penalty wrapper toggled off while testing for a memory leak:
#   - enable_time_penalty (N2): TRUE  -> FALSE
#   - time_penalty_per_step (O2): -1  -> 0
# The editor's cursor/selection was left on M13 after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQN_Runs")
$ws.Activate()

# enable_time_penalty -> FALSE
$ws.Range("N2").Value = $false

# time_penalty_per_step -> 0
$ws.Range("O2").Value = 0

# Leave the selection where the author's cursor ended up
$ws.Range("M13").Select()
